# Revert "Add new notes":
#   - removes the "Branches:" paragraph, the blank paragraph after it,
#     and the "Merge:" paragraph that were added by that commit
#   - the _GoBack bookmark that had ended up on the "Merge:" paragraph
#     moves back to the start of the (now last) "Use the body to explain
#     what and why vs. how" paragraph.
$d = $word.ActiveDocument

# Locate the paragraph that will become the new last paragraph.
$useBodyRange = $d.Content
$useBodyRange.Find.Execute("Use the body to explain what and why vs. how") | Out-Null
$useBodyStart = $useBodyRange.Start

$useBodyPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $useBodyStart -and $p.Range.End -gt $useBodyStart) {
        $useBodyPara = $p
        break
    }
}

# Re-seat the _GoBack bookmark at the very start of that paragraph
# (before its run), matching where it sits once the trailing
# "Branches:" / blank / "Merge:" paragraphs are gone. Adding a bookmark
# with the same name relocates it (bookmark names are unique).
$insertPoint = $d.Range($useBodyPara.Range.Start, $useBodyPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $insertPoint) | Out-Null

# Find the span covering "Branches:" through the end of "Merge:" and
# delete it outright (this removes the blank paragraph between them too).
$branchesRange = $d.Content
$branchesRange.Find.Execute("Branches:") | Out-Null
$branchesStart = $branchesRange.Start

$mergeRange = $d.Content
$mergeRange.Find.Execute("Merge:") | Out-Null
$mergeEnd = $mergeRange.End

$mergePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $mergeEnd -and $p.Range.End -ge $mergeEnd) {
        $mergePara = $p
        break
    }
}

$deleteRange = $d.Range($branchesStart, $mergePara.Range.End)
$deleteRange.Delete()
